$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'323.68"
$ws.Range("E2").Value = "'-2.53%"
$ws.Range("D3").Value = "'42.43"
$ws.Range("E3").Value = "'-7.24%"
$ws.Range("D4").Value = "'5.271"
$ws.Range("E4").Value = "'-7.40%"
$ws.Range("D5").Value = "'0.08144"
$ws.Range("E5").Value = "'-2.63%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.358"
$ws.Range("E6").Value = "'-2.53%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.815"
$ws.Range("E7").Value = "'-10.88%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9587"
$ws.Range("E8").Value = "'-2.42%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1127"
$ws.Range("E9").Value = "'-3.47%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1859"
$ws.Range("E10").Value = "'-4.40%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09321"
$ws.Range("E11").Value = "'-7.40%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04638"
$ws.Range("E12").Value = "'-0.53%"
$ws.Range("B13").Value = "MCDex"
$ws.Range("C13").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D13").Value = "'7.465"
$ws.Range("E13").Value = "'-27.99%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1064"
$ws.Range("E14").Value = "'0.55%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001291"
$ws.Range("E15").Value = "'0.58%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005874"
$ws.Range("E16").Value = "'-2.83%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.377"
$ws.Range("E17").Value = "'0.21%"
$ws.Range("D19").Value = "'0.3363"
$ws.Range("D20").Value = "'0.1375"
$ws.Range("E20").Value = "'-1.75%"
$ws.Range("D21").Value = "'0.2600"
$ws.Range("E21").Value = "'0.29%"
$ws.Range("D22").Value = "'0.04191"
$ws.Range("E22").Value = "'-0.29%"
$ws.Range("D23").Value = "'0.001259"
$ws.Range("E23").Value = "'-3.62%"
$ws.Range("D24").Value = "'0.004313"
$ws.Range("E24").Value = "'-5.95%"
$ws.Range("D25").Value = "'0.0001307"
$ws.Range("E25").Value = "'2.03%"
$ws.Range("D26").Value = "'0.0003012"
$ws.Range("E26").Value = "'-19.49%"
$ws.Range("D38").Value = "'0.02585"
$ws.Range("E38").Value = "'-7.00%"
$ws.Range("D39").Value = "'0.05470"
$ws.Range("E39").Value = "'-6.12%"
$ws.Range("D40").Value = "'0.007860"
$ws.Range("E40").Value = "'1.58%"
$ws.Range("D41").Value = "'0.1395"
$ws.Range("E41").Value = "'-2.86%"
$ws.Range("D42").Value = "'0.006632"
$ws.Range("E42").Value = "'-7.80%"
$ws.Range("D43").Value = "'0.002128"
$ws.Range("E43").Value = "'7.80%"
$ws.Range("D44").Value = "'0.008711"
$ws.Range("E44").Value = "'6.63%"
$ws.Range("D45").Value = "'0.3293"
$ws.Range("D46").Value = "'0.00007048"
$ws.Range("E46").Value = "'-2.32%"
$ws.Range("D47").Value = "'0.00000000758"
$ws.Range("E47").Value = "'1.08%"
$ws.Range("D48").Value = "'0.003502"
$ws.Range("E48").Value = "'0.26%"
$ws.Range("D49").Value = "'0.003569"
$ws.Range("E49").Value = "'2.01%"
$ws.Range("D50").Value = "'0.00002123"
$ws.Range("E50").Value = "'1.08%"
$ws.Range("D51").Value = "'0.0002022"
$ws.Range("E51").Value = "'1.08%"
